$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (English National League: Forest Green vs Tamworth FC)
$ws.Rows.Item(3).Insert()

# B3/C3 look like dates/times; force text format so Excel does not convert them to date/time serials
$ws.Range("B3:C3").NumberFormat = "@"

# Populate new row 3 with full data
$ws.Range("A3").Value = 'English National League'
$ws.Range("B3").Value = '2025-11-11'
$ws.Range("C3").Value = '16:45:00'
$ws.Range("D3").Value = 'Forest Green'
$ws.Range("E3").Value = 'Tamworth FC'
$ws.Range("F3").Value = 1.62
$ws.Range("G3").Value = 1.65
$ws.Range("H3").Value = 6
$ws.Range("I3").Value = 7.2
$ws.Range("J3").Value = 4.1
$ws.Range("K3").Value = 4.6
$ws.Range("L3").Value = 1.01
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 4.1
$ws.Range("O3").Value = 1.26
$ws.Range("P3").Value = 2.06
$ws.Range("Q3").Value = 1.75
$ws.Range("R3").Value = 1.42
$ws.Range("S3").Value = 2.9
$ws.Range("T3").Value = 1.84
$ws.Range("U3").Value = 1.98
$ws.Range("V3").Value = 1.16
$ws.Range("W3").Value = 2.5
$ws.Range("X3").Value = 19
$ws.Range("Y3").Value = 24
$ws.Range("Z3").Value = 55
$ws.Range("AA3").Value = 210
$ws.Range("AB3").Value = 9.4
$ws.Range("AC3").Value = 10.5
$ws.Range("AD3").Value = 26
$ws.Range("AE3").Value = 95
$ws.Range("AF3").Value = 10.5
$ws.Range("AG3").Value = 10.5
$ws.Range("AH3").Value = 980
$ws.Range("AI3").Value = 85
$ws.Range("AJ3").Value = 16
$ws.Range("AK3").Value = 17.5
$ws.Range("AL3").Value = 980
$ws.Range("AM3").Value = 140
$ws.Range("AN3").Value = 8.8
$ws.Range("AO3").Value = 130

# Apply the updated odds for the rows that shifted down (and row 2, which only had its own odds refreshed)
# Row 2 (was row 2): 4 updated values
$ws.Range("F2").Value = 1.63  # was 1.64
$ws.Range("J2").Value = 3.5  # was 3.2
$ws.Range("V2").Value = 1.18  # was 1.19
$ws.Range("AL2").Value = 42  # was 46

# Row 4 (was row 3): 14 updated values
$ws.Range("F4").Value = 1.99  # was 2.02
$ws.Range("I4").Value = 3.75  # was 3.8
$ws.Range("J4").Value = 4.1  # was 4
$ws.Range("R4").Value = 1.45  # was 1.46
$ws.Range("S4").Value = 2.8  # was 2.78
$ws.Range("V4").Value = 1.36  # was 1.35
$ws.Range("X4").Value = 21  # was 23
$ws.Range("Y4").Value = 20  # was 21
$ws.Range("Z4").Value = 1000  # was 32
$ws.Range("AB4").Value = 13.5  # was 14
$ws.Range("AC4").Value = 10.5  # was 11
$ws.Range("AD4").Value = 18.5  # was 19
$ws.Range("AF4").Value = 17  # was 17.5
$ws.Range("AI4").Value = 1000  # was 55

# Row 5 (was row 4): 11 updated values
$ws.Range("F5").Value = 2.42  # was 2.24
$ws.Range("G5").Value = 2.78  # was 3
$ws.Range("H5").Value = 3.25  # was 3.15
$ws.Range("K5").Value = 3.3  # was 3.75
$ws.Range("M5").Value = 1.1  # was 1.09
$ws.Range("N5").Value = 2.44  # was 2.66
$ws.Range("Q5").Value = 2.42  # was 2.4
$ws.Range("S5").Value = 4.3  # was 4.8
$ws.Range("T5").Value = 1.96  # was 1.98
$ws.Range("V5").Value = 1.32  # was 1.31
$ws.Range("W5").Value = 1.56  # was 1.51

# Row 6 (was row 5): 15 updated values
$ws.Range("F6").Value = 2.26  # was 2.28
$ws.Range("G6").Value = 2.74  # was 2.88
$ws.Range("H6").Value = 3.45  # was 3.4
$ws.Range("I6").Value = 4.6  # was 4.5
$ws.Range("J6").Value = 2.88  # was 2.62
$ws.Range("K6").Value = 3.6  # was 3.55
$ws.Range("L6").Value = 1.48  # was 1.49
$ws.Range("N6").Value = 2.5  # was 2.48
$ws.Range("O6").Value = 1.53  # was 1.54
$ws.Range("P6").Value = 1.5  # was 1.49
$ws.Range("Q6").Value = 2.36  # was 2.38
$ws.Range("R6").Value = 1.18  # was 1.17
$ws.Range("T6").Value = 2.08  # was 1.94
$ws.Range("U6").Value = 1.74  # was 1.73
$ws.Range("W6").Value = 1.61  # was 1.6

# Row 7 (was row 6): 2 updated values
$ws.Range("F7").Value = 3.55  # was 3.25
$ws.Range("G7").Value = 4.4  # was 4.8

